$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update predicted goals (columns C & D) with new LinearRegression-based predictions,
# and refresh the dependent comparison columns (G, H, I, J) plus one text correction (A32).

$ws.Range("C2").Value = 4
$ws.Range("G2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 2
$ws.Range("G3").Value = -1
$ws.Range("J3").Value = 1
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 4
$ws.Range("G7").Value = -1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("B8").Value = "England"
$ws.Range("C8").Value = 0
$ws.Range("B9").Value = "Ukraine"
$ws.Range("D9").Value = 5
$ws.Range("G9").Value = -1
$ws.Range("J9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("G11").Value = -1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("D12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("B13").Value = "Czechia"
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 2
$ws.Range("I13").Value = 0
$ws.Range("C14").Value = 2
$ws.Range("G14").Value = 1
$ws.Range("J14").Value = 0
$ws.Range("C15").Value = 4
$ws.Range("G15").Value = 1
$ws.Range("J15").Value = 1
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("G16").Value = -1
$ws.Range("C17").Value = 1
$ws.Range("G17").Value = 0
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 1
$ws.Range("B18").Value = "England"
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("G18").Value = 1
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 2
$ws.Range("B20").Value = "Ukraine"
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 2
$ws.Range("G20").Value = -1
$ws.Range("J20").Value = 1
$ws.Range("D21").Value = 0
$ws.Range("G21").Value = 1
$ws.Range("J21").Value = 0
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("G22").Value = -1
$ws.Range("J22").Value = 0
$ws.Range("B23").Value = "Czechia"
$ws.Range("D23").Value = 0
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("H24").Value = -1
$ws.Range("J24").Value = 0
$ws.Range("D25").Value = 4
$ws.Range("G25").Value = -1
$ws.Range("H25").Value = 1
$ws.Range("J25").Value = 0
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 0
$ws.Range("G30").Value = 1
$ws.Range("J30").Value = 0
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("J31").Value = 1
$ws.Range("A32").Value = "Denmark"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("A33").Value = "England"
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 1
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 1
$ws.Range("A35").Value = "Ukraine"
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 1
$ws.Range("G35").Value = -1
$ws.Range("J35").Value = 0
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 4
$ws.Range("G36").Value = -1
$ws.Range("J36").Value = 0
$ws.Range("A37").Value = "Czechia"
$ws.Range("C37").Value = 3
$ws.Range("D37").Value = 1
$ws.Range("G37").Value = 1
$ws.Range("J37").Value = 0
